$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.262.70'
$ws.Range("E2").Value = '  +2.59%  '
$ws.Range("D3").Value = '1.585.67'
$ws.Range("E3").Value = '  +1.43%  '
$ws.Range("E4").Value = '  +1.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.60'
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("E7").Value = '  +1.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.94'
$ws.Range("E8").Value = '  +6.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.251'
$ws.Range("E9").Value = '  +0.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0598'
$ws.Range("E10").Value = '  +0.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0887'
$ws.Range("E11").Value = '  +2.29%  '
$ws.Range("D12").Value = '1.812.74'
$ws.Range("E12").Value = '  +1.47%  '
$ws.Range("D13").Value = '1.596.67'
$ws.Range("E13").Value = '  +2.24%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.529'
$ws.Range("E14").Value = '  +1.74%  '
$ws.Range("E15").Value = '  -0.38%  '
$ws.Range("D16").Value = '28.269.24'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.15'
$ws.Range("E17").Value = '  +1.08%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '227.37'
$ws.Range("E18").Value = '  +1.06%  '
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.45'
$ws.Range("E20").Value = '  -0.97%  '
$ws.Range("E21").Value = '  +1.22%  '
$ws.Range("E22").Value = '  -1.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.32'
$ws.Range("E23").Value = '  -0.84%  '
$ws.Range("E24").Value = '  +0.58%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.06'
$ws.Range("E25").Value = '  +1.36%  '
$ws.Range("E26").Value = '  -0.15%  '
$ws.Range("E27").Value = '  -0.91%  '
$ws.Range("E28").Value = '  -1.24%  '
$ws.Range("E29").Value = '  +1.26%  '
$ws.Range("E30").Value = '  -0.19%  '
$ws.Range("E31").Value = '  +0.34%  '
$ws.Range("E32").Value = '  -0.09%  '
$ws.Range("E33").Value = '  -0.32%  '
$ws.Range("D34").Value = '1.398.02'
$ws.Range("E34").Value = '  -3.81%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.59'
$ws.Range("E35").Value = '  -1.43%  '
$ws.Range("E36").Value = '  -8.30%  '
$ws.Range("E37").Value = '  +2.02%  '
$ws.Range("E38").Value = '  -0.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.54'
$ws.Range("E39").Value = '  +9.05%  '
$ws.Range("E40").Value = '  -0.25%  '
$ws.Range("E41").Value = '  -0.78%  '
$ws.Range("E42").Value = '  +1.22%  '
$ws.Range("E43").Value = '  +1.07%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.59'
$ws.Range("E44").Value = '  -2.08%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.980'
$ws.Range("E45").Value = '  +1.23%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.18'
$ws.Range("E46").Value = '  -1.03%  '
$ws.Range("D47").Value = '1.722.55'
$ws.Range("E48").Value = '  +1.76%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '86.77'
$ws.Range("E49").Value = '  +0.42%  '
$ws.Range("E50").Value = '  +5.97%  '
$ws.Range("E51").Value = '  -0.75%  '
